# Update countries & provincias Spain
# Applies reordering of a few countries (their updated stats now rank them
# higher than neighbouring countries) plus small stat corrections, and
# bumps the "last updated" timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Group 1: Estonia now ranks above Egipto ---------------------------
# Row 59 becomes Estonia with refreshed stats, Egipto's (unchanged) stats
# shift down to row 60.
$ws.Range("A59").Value = "Estonia"
$ws.Range("B59").Value = 1097
$ws.Range("C59").Value = 58
$ws.Range("D59").Value = 62
$ws.Range("E59").Value = 1020
$ws.Range("F59").Value = 17
$ws.Range("G59").Value = 2
$ws.Range("H59").Value = 15

$ws.Range("A60").Value = "Egipto"
$ws.Range("B60").Value = 1070
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 241
$ws.Range("E60").Value = 758
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 71

# --- Minor stat-only updates -------------------------------------------
# Row 68: Moldavia
$ws.Range("D68").Value = 30
$ws.Range("E68").Value = 708

# Row 72: Bosnia y Herzegovina
$ws.Range("B72").Value = 632
$ws.Range("C72").Value = 8
$ws.Range("E72").Value = 581

# --- Group 2: Libano now ranks above Bulgaria and Azerbaiyan ------------
# Row 77 becomes Libano with refreshed stats, Bulgaria's (unchanged) stats
# shift to row 78, and Azerbaiyan's (unchanged) stats shift to row 79.
$ws.Range("A77").Value = "Libano"
$ws.Range("B77").Value = 527
$ws.Range("C77").Value = 7
$ws.Range("D77").Value = 54
$ws.Range("E77").Value = 455
$ws.Range("F77").Value = 26
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 18

$ws.Range("A78").Value = "Bulgaria"
$ws.Range("B78").Value = 522
$ws.Range("C78").Value = 19
$ws.Range("D78").Value = 37
$ws.Range("E78").Value = 467
$ws.Range("F78").Value = 26
$ws.Range("G78").Value = 1
$ws.Range("H78").Value = 18

$ws.Range("A79").Value = "Azerbaiyan"
$ws.Range("B79").Value = 521
$ws.Range("C79").Value = 0
$ws.Range("D79").Value = 32
$ws.Range("E79").Value = 484
$ws.Range("F79").Value = 17
$ws.Range("G79").Value = 0
$ws.Range("H79").Value = 5

# --- Group 3: Estado de Palestina now ranks above Senegal ---------------
# Row 102 becomes Estado de Palestina with refreshed stats, Senegal's
# (unchanged) stats shift down to row 103.
$ws.Range("A102").Value = "Estado de Palestina"
$ws.Range("B102").Value = 226
$ws.Range("C102").Value = 9
$ws.Range("D102").Value = 23
$ws.Range("E102").Value = 202
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 1

$ws.Range("A103").Value = "Senegal"
$ws.Range("B103").Value = 219
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 72
$ws.Range("E103").Value = 145
$ws.Range("F103").Value = 1
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 2

# --- Group 4: Banglades now ranks above Liechtenstein and Madagascar ----
# Row 128 becomes Banglades with refreshed stats, Liechtenstein's
# (unchanged) stats shift to row 129, and Madagascar's (unchanged) stats
# shift to row 130.
$ws.Range("A128").Value = "Banglades"
$ws.Range("B128").Value = 88
$ws.Range("C128").Value = 18
$ws.Range("D128").Value = 30
$ws.Range("E128").Value = 49
$ws.Range("F128").Value = 1
$ws.Range("G128").Value = 1
$ws.Range("H128").Value = 9

$ws.Range("A129").Value = "Liechtenstein"
$ws.Range("B129").Value = 77
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 0
$ws.Range("E129").Value = 76
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 1

$ws.Range("A130").Value = "Madagascar"
$ws.Range("B130").Value = 70
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 0
$ws.Range("E130").Value = 70
$ws.Range("F130").Value = 6
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 0

# --- Update "last refreshed" timestamp -----------------------------------
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 10:52"
